# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (want-to-go count) and a few "最低票价" (lowest price)
# values across the 展览 / 演出 / 本地生活 sheets, plus the 全部类型
# roll-up sheet that mirrors the same rows.

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 3005
$ws.Range("F7").Value = 248
$ws.Range("F8").Value = 16
$ws.Range("F10").Value = 7076
$ws.Range("F12").Value = 127
$ws.Range("F13").Value = 397
$ws.Range("F14").Value = 632
$ws.Range("F15").Value = 1542
$ws.Range("F16").Value = 2286
$ws.Range("G16").Value = 75.8
$ws.Range("F17").Value = 1547
$ws.Range("G17").Value = 78
$ws.Range("F20").Value = 159
$ws.Range("F21").Value = 12
$ws.Range("F23").Value = 366
$ws.Range("F24").Value = 56
$ws.Range("F25").Value = 56
$ws.Range("F26").Value = 1803
$ws.Range("F27").Value = 1725
$ws.Range("F30").Value = 1694
$ws.Range("F31").Value = 1265
$ws.Range("F32").Value = 150
$ws.Range("F34").Value = 20
$ws.Range("F36").Value = 455
$ws.Range("F37").Value = 40
$ws.Range("F38").Value = 2528
$ws.Range("F39").Value = 2797
$ws.Range("F40").Value = 2085
$ws.Range("F41").Value = 42
$ws.Range("F43").Value = 1141
$ws.Range("F44").Value = 1
$ws.Range("F46").Value = 33
$ws.Range("F47").Value = 339
$ws.Range("F49").Value = 184
$ws.Range("F50").Value = 420

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G6").Value = "不可售"
$ws.Range("G8").Value = "不可售"
$ws.Range("F10").Value = 189
$ws.Range("F13").Value = 72
$ws.Range("F18").Value = 72
$ws.Range("F21").Value = 491
$ws.Range("F22").Value = 48
$ws.Range("F27").Value = 8
$ws.Range("F32").Value = 14

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1735
$ws.Range("F8").Value = 2795
$ws.Range("F9").Value = 1061
$ws.Range("F10").Value = 975
$ws.Range("F12").Value = 338
$ws.Range("F13").Value = 1657
$ws.Range("F14").Value = 7620

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 3005
$ws.Range("F6").Value = 1735
$ws.Range("F7").Value = 2795
$ws.Range("F8").Value = 7076
$ws.Range("F9").Value = 1061
$ws.Range("F11").Value = 397
$ws.Range("F12").Value = 632
$ws.Range("F13").Value = 1542
$ws.Range("F14").Value = 2286
$ws.Range("G14").Value = 75.8
$ws.Range("F15").Value = 1547
$ws.Range("G15").Value = 78
$ws.Range("F17").Value = 189
$ws.Range("F18").Value = 159
$ws.Range("F19").Value = 12
$ws.Range("F21").Value = 56
$ws.Range("F22").Value = 56
$ws.Range("F23").Value = 1803
$ws.Range("F24").Value = 72
$ws.Range("F27").Value = 1694
$ws.Range("F28").Value = 1265
$ws.Range("F29").Value = 150
$ws.Range("F31").Value = 20
$ws.Range("F33").Value = 72
$ws.Range("F36").Value = 455
$ws.Range("F37").Value = 40
$ws.Range("F38").Value = 2528
$ws.Range("F39").Value = 2797
$ws.Range("F40").Value = 2085
$ws.Range("F41").Value = 42
$ws.Range("F43").Value = 1141
$ws.Range("F45").Value = 33
$ws.Range("F47").Value = 184
$ws.Range("F49").Value = 420
